$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of car reservation dialogue data.
# Set C7 first, then A7, so the new shared-string entries are appended
# in the same order as the target workbook (C7's text becomes the first
# new unique string, A7's text becomes the second).
$ws.Range("C7").Value = "great let me do the reservation"
$ws.Range("A7").Value = "I see. thank you."

# Widen column C to fit the new content
$ws.Columns.Item(3).ColumnWidth = 46.8333333

# Move active selection to the newly added row
[void]$ws.Range("A7").Select()
